$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row above row 42 (shifts rows 42:163 down to 43:164,
# pushing the previous last row of data into the new row 164).
$ws.Rows("42:42").Insert()

# Populate the newly inserted row 42 with the new weekly price record.
$ws.Range("A42").Value = 7
$ws.Range("B42").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C42").Value = "Ñuble"
$ws.Range("D42").Value = 44560
$ws.Range("E42").Value = 16
$ws.Range("F42").Value = 100112017
$ws.Range("G42").Value = "Apio"
$ws.Range("H42").Value = "Americana (o)"
$ws.Range("I42").Value = "Primera"
$ws.Range("J42").Value = 60
$ws.Range("K42").Value = 8000
$ws.Range("L42").Value = 8500
$ws.Range("M42").Value = 8250
$ws.Range("N42").Value = "`$/docena de matas"
$ws.Range("O42").Value = "Provincia del Elquí"
$ws.Range("P42").Value = 1375
$ws.Range("Q42").Value = 6
$ws.Range("R42").Value = "Hortaliza"
